# Generate Report for Handoff
# The "b.md" row moves from "Handed back: in sync with en-US" to
# "Ready for handoff" on all three sheets, gets a fresh "Latest Handoff
# File" / "Latest Handoff Datetime" pair, and (on the per-locale sheets)
# an Error Detail noting the handback file is stale.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c8b1c1dbdb4743181aedb6c0de7b2069f45e86ca/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7076bfeefee084e7c5a81f055d40228dc2282297/e2e/b.md."

# ---- Overview sheet: row 3 is b.md ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-30 04:37:17"

# ---- zh-cn sheet: row 3 is b.md ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-30 04:37:12"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666666

# ---- de-de sheet: row 3 is b.md ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-30 04:37:17"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666666
